$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three header rows (Sekretariat ..., address, phone) - content below shifts up.
$ws.Rows("1:3").Delete() | Out-Null

# Restore the original per-row heights (row height is independent of the
# row's content and was not shifted when the rows above were removed).
$ws.Rows("1").RowHeight = 27.95
$ws.Rows("2").RowHeight = 27
$ws.Rows("3").RowHeight = 18
$ws.Rows("4").RowHeight = 32.1
$ws.Range("5:12").RowHeight = 21.95
$ws.Rows("13").RowHeight = 21.95
$ws.Rows("14").RowHeight = 23.1
$ws.Rows("15").RowHeight = 21

# Reselect to match the saved selection state.
$ws.Range("A1:G12").Select() | Out-Null
